# Add a new "latest run" row at the top of the data table (row 2), pushing
# all existing model/real-value rows down by one, and refresh the selection.
#
# Commit: "try to add all centrality operation in networkx" — a new
# centrality-run result (17_07_06_01_01_31_0_1_500) is inserted as the first
# data row, ahead of the previously-first run (17_06_23_10_09_48_1_500_500).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new blank row above the current row 2 — this shifts rows 2:12
# down to 3:13 (formulas, number formats and fill styles move with them).
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new run's label and measurements.
$ws.Range("A2").Value = "17_07_06_01_01_31_0_1_500"
$ws.Range("B2").Value = 0.227858305734011
$ws.Range("C2").Value = 0.191372665244543
$ws.Range("D2").Value = 0.379406709533436
$ws.Range("E2").Value = 0.144618203673806
$ws.Range("F2").Value = 0.0492559229829956
$ws.Range("G2").Value = 0.0689063882558996
$ws.Range("H2").Value = 0.110164559551915
$ws.Range("I2").Value = 0.0921915925927093
$ws.Range("J2").Formula = "=SUM(B2:I2)"

# Match the author's final selection (cell J9 was last clicked).
$ws.Range("J9").Select() | Out-Null
